$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-08-24 Thursday" "2023-08-25 Friday"

Replace-Text "44×34=" "15×30="
Replace-Text "62×25=" "42×35="
Replace-Text "92×45=" "39×23="
Replace-Text "54×36=" "48×20="
Replace-Text "79×55=" "41×81="
Replace-Text "92×72=" "52×37="
Replace-Text "14×84=" "39×21="
Replace-Text "34×72=" "92×16="
Replace-Text "49×14=" "57×85="
Replace-Text "42×85=" "91×69="
Replace-Text "95×79=" "38×24="
Replace-Text "73×33=" "21×84="
Replace-Text "80×93=" "20×88="
Replace-Text "99×27=" "28×63="
Replace-Text "23×22=" "43×98="
Replace-Text "59×21=" "79×25="
Replace-Text "59×22=" "45×60="
Replace-Text "48×83=" "42×27="
Replace-Text "67×50=" "63×92="
Replace-Text "33×65=" "91×63="
Replace-Text "15×75=" "33×23="
Replace-Text "68×67=" "15×24="
Replace-Text "23×32=" "76×48="
Replace-Text "42×16=" "18×37="
Replace-Text "96×17=" "52×79="
